{"js": "// Update the date line and the 24 two-digit-by-two-digit multiplication\n// answers in the table to the new problems/answers from the commit.\nconst replacements = [\n  [\"2025-05-15 Thursday\", \"2025-05-16 Friday\"],\n  [\"22\\u00d747=1034\", \"45\\u00d743=1935\"],\n  [\"92\\u00d727=2484\", \"76\\u00d772=5472\"],\n  [\"74\\u00d749=3626\", \"86\\u00d759=5074\"],\n  [\"16\\u00d761=976\", \"48\\u00d773=3504\"],\n  [\"88\\u00d742=3696\", \"16\\u00d758=928\"],\n  [\"88\\u00d741=3608\", \"37\\u00d795=3515\"],\n  [\"98\\u00d737=3626\", \"38\\u00d751=1938\"],\n  [\"24\\u00d711=264\", \"18\\u00d784=1512\"],\n  [\"19\\u00d795=1805\", \"77\\u00d793=7161\"],\n  [\"56\\u00d771=3976\", \"49\\u00d773=3577\"],\n  [\"64\\u00d738=2432\", \"90\\u00d770=6300\"],\n  [\"63\\u00d796=6048\", \"73\\u00d715=1095\"],\n  [\"35\\u00d748=1680\", \"72\\u00d765=4680\"],\n  [\"78\\u00d762=4836\", \"45\\u00d791=4095\"],\n  [\"59\\u00d762=3658\", \"38\\u00d795=3610\"],\n  [\"66\\u00d727=1782\", \"80\\u00d767=5360\"],\n  [\"23\\u00d718=414\", \"76\\u00d797=7372\"],\n  [\"41\\u00d784=3444\", \"52\\u00d756=2912\"],\n  [\"31\\u00d712=372\", \"17\\u00d775=1275\"],\n  [\"15\\u00d731=465\", \"28\\u00d768=1904\"],\n  [\"39\\u00d790=3510\", \"38\\u00d769=2622\"],\n  [\"68\\u00d725=1700\", \"65\\u00d746=2990\"],\n  [\"11\\u00d716=176\", \"29\\u00d741=1189\"],\n  [\"78\\u00d729=2262\", \"53\\u00d723=1219\"],\n  [\"14\\u00d757=798\", \"11\\u00d776=836\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 24 two-digit-by-two-digit multiplication\n# answers in the table to the new problems/answers from the commit.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-05-15 Thursday\", \"2025-05-16 Friday\"),\n    @(\"22\u00d747=1034\", \"45\u00d743=1935\"),\n    @(\"92\u00d727=2484\", \"76\u00d772=5472\"),\n    @(\"74\u00d749=3626\", \"86\u00d759=5074\"),\n    @(\"16\u00d761=976\", \"48\u00d773=3504\"),\n    @(\"88\u00d742=3696\", \"16\u00d758=928\"),\n    @(\"88\u00d741=3608\", \"37\u00d795=3515\"),\n    @(\"98\u00d737=3626\", \"38\u00d751=1938\"),\n    @(\"24\u00d711=264\", \"18\u00d784=1512\"),\n    @(\"19\u00d795=1805\", \"77\u00d793=7161\"),\n    @(\"56\u00d771=3976\", \"49\u00d773=3577\"),\n    @(\"64\u00d738=2432\", \"90\u00d770=6300\"),\n    @(\"63\u00d796=6048\", \"73\u00d715=1095\"),\n    @(\"35\u00d748=1680\", \"72\u00d765=4680\"),\n    @(\"78\u00d762=4836\", \"45\u00d791=4095\"),\n    @(\"59\u00d762=3658\", \"38\u00d795=3610\"),\n    @(\"66\u00d727=1782\", \"80\u00d767=5360\"),\n    @(\"23\u00d718=414\", \"76\u00d797=7372\"),\n    @(\"41\u00d784=3444\", \"52\u00d756=2912\"),\n    @(\"31\u00d712=372\", \"17\u00d775=1275\"),\n    @(\"15\u00d731=465\", \"28\u00d768=1904\"),\n    @(\"39\u00d790=3510\", \"38\u00d769=2622\"),\n    @(\"68\u00d725=1700\", \"65\u00d746=2990\"),\n    @(\"11\u00d716=176\", \"29\u00d741=1189\"),\n    @(\"78\u00d729=2262\", \"53\u00d723=1219\"),\n    @(\"14\u00d757=798\", \"11\u00d776=836\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $true, $true, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
